# Applies the "Added a 11 phrases and tweaked 3" commit to BINGO_cc.xlsx
# - tweak 3 existing phrases in the "list" sheet
# - append 11 new phrases to the "list" sheet
# - move the cell selection on "original" (A:B2 -> E4) and "list" (A45 -> A46)

$wb = $excel.ActiveWorkbook

$wsOriginal = $wb.Worksheets.Item("original")
$wsList     = $wb.Worksheets.Item("list")

# --- tweak existing phrases -------------------------------------------------
$wsList.Range("A3").Value  = "My outlook / WebEx is not working"
$wsList.Range("A40").Value = "I" + [char]0x2019 + "ll take silence as _____"
$wsList.Range("A41").Value = "You're breaking up"
$wsList.Range("A45").Value = "I have / they had to drop"

# --- append 11 new phrases ---------------------------------------------------
$newPhrases = @(
    "Unnecessary verbing (""Let's solution that"")",
    "It's on my radar",
    "It is what it is",
    "Ping me",
    "[I don't] have the bandwidth",
    "Ducks in a row",
    "Chewing sounds",
    "Animal or child sounds",
    "You cut out",
    "It's loading",
    "Next slide please"
)

$row = 46
foreach ($phrase in $newPhrases) {
    $wsList.Range("A" + $row).Value = $phrase
    $row = $row + 1
}

# --- selection / view bookkeeping -------------------------------------------
# Move the "original" sheet's remembered selection without stealing focus
# away from the "list" sheet (which stays the active tab).
$wsOriginal.Range("E4").Select()
$wsList.Activate()
$wsList.Range("A46").Select()
